# Updates cryptos price/volume data (Price column D, Volume(1h) column E)
# for the rows whose values changed, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.880.86'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.85%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.213.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -6.83%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.21%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '79.69'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -9.46%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.79%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.456'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -7.58%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0772'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.78%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '27.83'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -10.07%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.66'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -14.02%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.96%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.545.19'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -7.03%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.09'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -7.36%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.89'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -7.51%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.242.75'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.81%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.710'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.79%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '38.756.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.84%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0857'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.00%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.71'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.75'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.63%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.77'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -9.38%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '225.17'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.16%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -10.13%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.90%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.17'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.83%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.09%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '147.55'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.96'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -8.31%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.21%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.74'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -9.01%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.31'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0686'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -6.41%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.64'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.93%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0946'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.44%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.43'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -9.42%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.58'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -8.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.57'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.58%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.905.51'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.55%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.04'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -13.47%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0253'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.49%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.07'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -8.28%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.85'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.68%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.51'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -9.28%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.421.89'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.96%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.04'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.81%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '87.01'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.13%  '

